# AddressingIssues.xlsx - "Changes after MSU-IIT work."
#
# Summary of the edit:
#  - Entry sheet (sheet2) gains a new column D with source-code comment
#    strings, a yellow highlight style on column B (and every 3rd C cell),
#    and becomes the active/selected sheet (selection moves to C3).
#  - Pages sheet (sheet1) selection moves to I15.
#  - Content Per Page sheet (sheet3) is no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsPages = $wb.Worksheets.Item(1)
$wsEntry = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# New column D content for the Entry sheet - ported-over C header/enum
# comments that document the byte values used by the Entry table.
# Kept as an ordered list (not a hashtable) so the shared-string table
# is built in exactly this row order.
# ---------------------------------------------------------------------
$dValues = @(
    @{ Row = 2;  Text = '        /*      unsigned char   ss ;    // seconds' }
    @{ Row = 3;  Text = '        unsigned char   mn ;    // minutes' }
    @{ Row = 4;  Text = '        unsigned char   hh ;    // hours' }
    @{ Row = 5;  Text = '        unsigned char   md ;    // day in month, from 1 to 31' }
    @{ Row = 6;  Text = '        unsigned char   wd ;    // day in week, monday=0, tuesday=1, .... sunday=6' }
    @{ Row = 7;  Text = '        unsigned char   mo ;    // month number, from 1 to 12 (and not from 0 to 11 as with unix C time !)' }
    @{ Row = 8;  Text = '        unsigned int    yy ;    // year Y2K compliant, from 1892 to 2038*/' }
    @{ Row = 10; Text = 'case 255: return "*";' }
    @{ Row = 11; Text = 'case 254: return "-";' }
    @{ Row = 12; Text = 'case 253: return "/";' }
    @{ Row = 13; Text = 'case 252: return "On";' }
    @{ Row = 14; Text = 'case 251: return "Off";' }
    @{ Row = 15; Text = 'case 250: return "";' }
    @{ Row = 17; Text = 'ASTERISK = 255,' }
    @{ Row = 18; Text = 'DASH = 254,' }
    @{ Row = 19; Text = 'SLASH = 253,' }
    @{ Row = 20; Text = 'ON = 252,' }
    @{ Row = 21; Text = 'OFF = 251,' }
    @{ Row = 22; Text = 'EMPTY = 250' }
)

foreach ($entry in $dValues) {
    $wsEntry.Range("D$($entry.Row)").Value = $entry.Text
}

# ---------------------------------------------------------------------
# Highlight column B (rows 2-22) with a yellow fill, and every third
# entry in column C (rows 6, 9, 12, 15, 18, 21) with the same fill.
# ---------------------------------------------------------------------
$wsEntry.Range("B2:B22").Interior.Color = 65535

$highlightedCRows = @(6, 9, 12, 15, 18, 21)
foreach ($row in $highlightedCRows) {
    $wsEntry.Range("C$row").Interior.Color = 65535
}

# ---------------------------------------------------------------------
# Selections / active sheet: Pages keeps its own cursor position, but
# the Entry sheet becomes the active tab with C3 selected.
# ---------------------------------------------------------------------
$wsPages.Range("I15").Select()
$wsEntry.Range("C3").Select()
